# repull data, push all data, mean calculation
# Update column F (dSF) values for the affected rows to reflect the repulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 1
$ws.Range("F3").Value = -1
$ws.Range("F4").Value = -7
$ws.Range("F5").Value = -3
$ws.Range("F6").Value = 1
$ws.Range("F7").Value = 2
$ws.Range("F8").Value = -1
$ws.Range("F9").Value = -2
$ws.Range("F10").Value = -4
$ws.Range("F11").Value = -2
$ws.Range("F12").Value = 1
$ws.Range("F13").Value = 2
$ws.Range("F14").Value = 0
$ws.Range("F16").Value = -3
$ws.Range("F17").Value = 5
$ws.Range("F18").Value = -5
$ws.Range("F21").Value = 7
$ws.Range("F22").Value = -8
$ws.Range("F23").Value = 6
$ws.Range("F27").Value = -1
$ws.Range("F28").Value = 2
$ws.Range("F31").Value = -1
$ws.Range("F33").Value = 0
$ws.Range("F34").Value = 2
$ws.Range("F36").Value = 4
$ws.Range("F39").Value = 1
$ws.Range("F40").Value = 4
$ws.Range("F42").Value = 4
$ws.Range("F43").Value = 0
